# Apply crypto price/volume updates for Sat Jun  1 22:48:14 UTC 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.781.27"
$ws.Range("E2").Value = "  +0.24%  "

# Row 3
$ws.Range("D3").Value = "3.825.16"
$ws.Range("E3").Value = "  +1.18%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.93"
$ws.Range("D5").ClearFormats()

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.79"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.25%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("E8").Value = "  +0.13%  "

# Row 9
$ws.Range("E9").Value = "  +0.65%  "

# Row 10
$ws.Range("E10").Value = "  -0.91%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.452"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.81%  "

# Row 12
$ws.Range("E12").Value = "  -0.53%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.94"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.83%  "

# Row 14
$ws.Range("D14").Value = "4.469.93"
$ws.Range("E14").Value = "  +1.27%  "

# Row 15
$ws.Range("D15").Value = "3.784.16"
$ws.Range("E15").Value = "  +0.48%  "

# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "67.827.74"
$ws.Range("E16").Value = "  +0.36%  "

# Row 17
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.47"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.13%  "

# Row 18
$ws.Range("E18").Value = "  +1.41%  "

# Row 19
$ws.Range("E19").Value = "  +0.60%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.34"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.49%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.94"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.54%  "

# Row 22
$ws.Range("E22").Value = "  +0.52%  "

# Row 23
$ws.Range("E23").Value = "  -3.28%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.44"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.11%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.09"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.52%  "

# Row 26
$ws.Range("E26").Value = "  -0.89%  "

# Row 27
$ws.Range("E27").Value = "  +0.35%  "

# Row 28
$ws.Range("E28").Value = "  +0.21%  "

# Row 29
$ws.Range("D29").Value = "3.973.52"

# Row 30
$ws.Range("E30").Value = "  +0.13%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.40"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.81%  "

# Row 32
$ws.Range("E32").Value = "  +1.78%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.72"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.16%  "

# Row 34
$ws.Range("E34").Value = "  -0.05%  "

# Row 35
$ws.Range("B35").Value = "RenzoRestakedETH"
$ws.Range("C35").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D35").Value = "3.764.84"
$ws.Range("E35").Value = "  +0.83%  "

# Row 36
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.09"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.34%  "

# Row 37
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.100"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.02%  "

# Row 38
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.36"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.14%  "

# Row 39
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.138"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.44%  "

# Row 40
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.71%  "

# Row 41
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.81"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.98%  "

# Row 42
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.04%  "

# Row 43
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.01%  "

# Row 44
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.10"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.10%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "28.81"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +12.67%  "

# Row 46
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.301"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.45%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.41"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +12.79%  "

# Row 48
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "43.13"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.95%  "

# Row 49
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.36"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.23%  "

# Row 50
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "148.23"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.12%  "

# Row 51
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.84"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.38%  "

